$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 25 (the row about to be pushed down) so formatting (incl. borders)
# is preserved exactly for both rows once we insert a new one.
$ws.Rows.Item(25).Copy()
$ws.Rows.Item(25).Insert()

# Fill in the new item's data in row 25
$ws.Cells.Item(25, 1).Value = 19
$ws.Cells.Item(25, 3).Value = "D.DEP 10.000 I.U. 30 CAPS"
$ws.Cells.Item(25, 8).Value = "0:1"
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 14).Value = "300.00"
$ws.Cells.Item(25, 16).Value = "150.0000"
$ws.Cells.Item(25, 17).Value = "0:1"

Write-Host "Done inserting row"
